# Update the cryptos price/volume table with the latest scraped values.
# Values are written as text (not numbers) so strings such as "26.844.38"
# or "1.000" keep their exact original formatting instead of being
# coerced into numeric values by Excel's automatic type detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "26.844.38"
Set-TextValue "E2" "  -1.39%  "
Set-TextValue "D3" "1.874.17"
Set-TextValue "E3" "  -1.61%  "
Set-TextValue "E4" "  -0.24%  "
Set-TextValue "D5" "301.14"
Set-TextValue "E5" "  -2.17%  "
Set-TextValue "D7" "0.5346"
Set-TextValue "E7" "  +1.79%  "
Set-TextValue "E8" "  -1.65%  "
Set-TextValue "D9" "0.07179"
Set-TextValue "E9" "  -1.77%  "
Set-TextValue "D10" "21.65"
Set-TextValue "D11" "0.8881"
Set-TextValue "E11" "  -1.95%  "
Set-TextValue "D12" "0.08151"
Set-TextValue "E12" "  +0.90%  "
Set-TextValue "D13" "1.910.35"
Set-TextValue "E13" "  +4.14%  "
Set-TextValue "D14" "93.31"
Set-TextValue "E14" "  -2.73%  "
Set-TextValue "D15" "5.298"
Set-TextValue "E15" "  -1.36%  "
Set-TextValue "E16" "  -0.21%  "
Set-TextValue "E17" "  +0.32%  "
Set-TextValue "E18" "  -1.54%  "
Set-TextValue "D20" "26.885.88"
Set-TextValue "E20" "  -1.36%  "
Set-TextValue "D22" "10.66"
Set-TextValue "E22" "  -1.50%  "
Set-TextValue "D23" "6.400"
Set-TextValue "E23" "  -1.37%  "
Set-TextValue "D24" "146.37"
Set-TextValue "E24" "  -2.33%  "
Set-TextValue "D25" "2.280"
Set-TextValue "E25" "  -3.37%  "
Set-TextValue "B26" "EthereumClassic"
Set-TextValue "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "18.05"
Set-TextValue "E26" "  -1.17%  "
Set-TextValue "B27" "Toncoin"
Set-TextValue "C27" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D27" "1.728"
Set-TextValue "E27" "  -0.91%  "
Set-TextValue "D28" "113.98"
Set-TextValue "E28" "  -2.59%  "
Set-TextValue "E29" "  -2.47%  "
Set-TextValue "D30" "4.616"
Set-TextValue "E30" "  -5.49%  "
Set-TextValue "D31" "0.09135"
Set-TextValue "E31" "  -1.15%  "
Set-TextValue "D32" "0.8163"
Set-TextValue "E32" "  +1.21%  "
Set-TextValue "D33" "0.04988"
Set-TextValue "E33" "  -1.61%  "
Set-TextValue "D34" "1.175"
Set-TextValue "E34" "  -4.51%  "
Set-TextValue "D35" "2.950"
Set-TextValue "E35" "  -0.96%  "
Set-TextValue "D36" "0.6071"
Set-TextValue "E36" "  +5.78%  "
Set-TextValue "D37" "3.193"
Set-TextValue "E37" "  -5.79%  "
Set-TextValue "D38" "2.614"
Set-TextValue "E38" "  -3.56%  "
Set-TextValue "E39" "  -2.17%  "
Set-TextValue "D40" "1.070"
Set-TextValue "E40" "  -1.50%  "
Set-TextValue "D41" "6.619"
Set-TextValue "E41" "  +0.14%  "
Set-TextValue "D42" "8.934"
Set-TextValue "E42" "  -0.78%  "
Set-TextValue "D43" "0.5149"
Set-TextValue "E43" "  +4.81%  "
Set-TextValue "D44" "114.98"
Set-TextValue "E44" "  -1.33%  "
Set-TextValue "D45" "0.1493"
Set-TextValue "E45" "  -1.54%  "
Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  -0.25%  "
Set-TextValue "E47" "  -0.69%  "
Set-TextValue "D48" "9.906"
Set-TextValue "E48" "  -2.96%  "
Set-TextValue "D49" "37.53"
Set-TextValue "E49" "  -2.63%  "
Set-TextValue "D50" "0.06068"
Set-TextValue "E50" "  +1.75%  "
Set-TextValue "D51" "62.22"
